$d = $word.ActiveDocument

# 1. Update the NOTE paragraph: change the async version note text and
#    remove the "2.0" wording in favor of the simpler 2.1.x note.
$old1 = " in the solutions folder is 1.4.x and the version obtained with npm install async is version 2.0.0-rc3.  These are NOT compatible.  Feel free to examine the online documentation to determine the differences.  HINT: the order of the parameters is different.  We will change the labs when 2.0 is released.  If the task has NO dependencies, it has only one parameter, the callback.  If the task has dependencies, the task has two parameters, the current results and the callback."
$new1 = " in the solutions folder is 2.1.x.  If the task has NO dependencies, it has only one parameter, the callback.  If the task has dependencies, the task has two parameters, the current results first, then the callback."

$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 2. Step 5 paragraph: "your exists" -> "the exists"
$old2 = "This time, you should observe your exists function being invoked"
$new2 = "This time, you should observe the exists function being invoked"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# 3. Remove the _GoBack bookmark from the NOTE paragraph (it's no longer there)
#    and add a new one inside the "minification" sentence, right after the word
#    "minification".
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

$findRange = $d.Content
$findRange.Find.Execute("minification workflow by adding two more tasks")
$bookmarkPos = $findRange.Start + 12
$markRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $markRange)
